$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. "Issue Contributors" heading -> bump font size to 14pt (w:sz 28)
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Issue Contributors", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $rng.Font.Size = 14
}

# ---------------------------------------------------------------
# 2. Description paragraph -> append a new sentence
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("move along its life cycle.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" It is important for these contributors to be organized in a coherent and logical fashion. Effective division of labor among project members can drastically improve the project efficiency.")
}

# ---------------------------------------------------------------
# 3. Remove the stray "_GoBack" bookmark (Triggers paragraph)
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 4. Actors paragraph -> append a trailing period
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Project managers and users working on the project", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $rng.Collapse(0)
    $rng.InsertAfter(".")
}

# ---------------------------------------------------------------
# 5. Preconditions paragraph -> fill in "None. " (was empty)
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Preconditions", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $headingPara = $rng.Paragraphs.First
    $answerPara = $headingPara.Next
    $answerPara.Range.InsertBefore("None. ")
}

# ---------------------------------------------------------------
# 6. Main Success Scenario paragraph -> rewrite ending
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("the projects issue.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "the project issues, experience levels, and interest.", 2) | Out-Null
